$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.8367010369953505
$ws.Range("C2").Value = 0.2049245934894088
$ws.Range("D2").Value = 0.3510291671064465
$ws.Range("F2").Value = 0.8901291982779966
$ws.Range("G2").Value = 0.3445974802329985
$ws.Range("H2").Value = 0.5178742945755914
$ws.Range("I2").Value = 0.360598179036014
$ws.Range("J2").Value = 0.2776838763725493
$ws.Range("M2").Value = 0.4061749507860384
$ws.Range("O2").Value = 1.660588337173365
# Row 3
$ws.Range("B3").Value = 0.7328516411075725
$ws.Range("C3").Value = 0.1796925741805921
$ws.Range("D3").Value = 0.3476827487559291
$ws.Range("F3").Value = 0.8954113168092306
$ws.Range("G3").Value = 0.3475926886938012
$ws.Range("H3").Value = 0.5233507373857194
$ws.Range("I3").Value = 0.3684798755134917
$ws.Range("J3").Value = 0.2794330165164993
$ws.Range("M3").Value = 0.3743556449220549
$ws.Range("O3").Value = 1.678139633395617
# Row 4
$ws.Range("B4").Value = 0.6688563402614989
$ws.Range("C4").Value = 0.1641304565277721
$ws.Range("D4").Value = 0.3457661527478137
$ws.Range("F4").Value = 0.8992921972119419
$ws.Range("G4").Value = 0.3497855927703668
$ws.Range("H4").Value = 0.5270134162765032
$ws.Range("I4").Value = 0.3736292901093865
$ws.Range("J4").Value = 0.280683080905348
$ws.Range("M4").Value = 0.3548479801024413
$ws.Range("O4").Value = 1.690285523443805
# Row 5
$ws.Range("B5").Value = 0.6427214702790991
$ws.Range("C5").Value = 0.1577717523015281
$ws.Range("D5").Value = 0.3450199661950251
$ws.Range("F5").Value = 0.901033914261447
$ws.Range("G5").Value = 0.3507679744042562
$ws.Range("H5").Value = 0.528581437393477
$ws.Range("I5").Value = 0.3758055347362106
$ws.Range("J5").Value = 0.2812367387887136
$ws.Range("M5").Value = 0.3469064176336261
$ws.Range("O5").Value = 1.69557889150245
# Row 6
$ws.Range("B6").Value = 0.6383784532146421
$ws.Range("C6").Value = 0.1567148807035323
$ws.Range("D6").Value = 0.3448981701171476
$ws.Range("F6").Value = 0.9013327990403965
$ws.Range("G6").Value = 0.3509364522931264
$ws.Range("H6").Value = 0.5288463625507447
$ws.Range("I6").Value = 0.3761715942756432
$ws.Range("J6").Value = 0.2813313447127896
$ws.Range("M6").Value = 0.3455882247156339
$ws.Range("O6").Value = 1.696478602124998
# Row 7
$ws.Range("B7").Value = 0.6685041015389857
$ws.Range("C7").Value = 0.1640447690785152
$ws.Range("D7").Value = 0.3457559481941956
$ws.Range("F7").Value = 0.8993150380153949
$ws.Range("G7").Value = 0.3497984824369311
$ws.Range("H7").Value = 0.5270342576895004
$ws.Range("I7").Value = 0.3736583248171588
$ws.Range("J7").Value = 0.2806903686181741
$ws.Range("M7").Value = 0.3547408444189912
$ws.Range("O7").Value = 1.690355520222326
# Row 8
$ws.Range("B8").Value = 0.8009428929544242
$ws.Range("C8").Value = 0.1962393062387946
$ws.Range("D8").Value = 0.3498467203103246
$ws.Range("F8").Value = 0.8918180621691008
$ws.Range("G8").Value = 0.3455566580474496
$ws.Range("H8").Value = 0.5197002712655348
$ws.Range("I8").Value = 0.3632513594241633
$ws.Range("J8").Value = 0.2782504141242796
$ws.Range("M8").Value = 0.3951978449998705
$ws.Range("O8").Value = 1.666355515581245
# Row 9
$ws.Range("B9").Value = 1.058747988434163
$ws.Range("C9").Value = 0.2588029451156331
$ws.Range("D9").Value = 0.3589604801805564
$ws.Range("F9").Value = 0.882181705006424
$ws.Range("G9").Value = 0.3400560806222046
$ws.Range("H9").Value = 0.5077005186676899
$ws.Range("I9").Value = 0.3453092140955132
$ws.Range("J9").Value = 0.2748644334075507
$ws.Range("M9").Value = 0.4747481226084815
$ws.Range("O9").Value = 1.630180035084123
# Row 10
$ws.Range("B10").Value = 1.24691758343107
$ws.Range("C10").Value = 0.304401573667775
$ws.Range("D10").Value = 0.3663174476890418
$ws.Range("F10").Value = 0.8781988678379236
$ws.Range("G10").Value = 0.3377468416552745
$ws.Range("H10").Value = 0.5003377588781746
$ws.Range("I10").Value = 0.3336381868611706
$ws.Range("J10").Value = 0.2732319598944102
$ws.Range("M10").Value = 0.5333033945194501
$ws.Range("O10").Value = 1.61027313543849
# Row 11
$ws.Range("B11").Value = 1.332236530894306
$ws.Range("C11").Value = 0.3250620813572027
$ws.Range("D11").Value = 0.3698069911630597
$ws.Range("F11").Value = 0.8770614689680798
$ws.Range("G11").Value = 0.3370754832465224
$ws.Range("H11").Value = 0.4973041258549742
$ws.Range("I11").Value = 0.3286584960493926
$ws.Range("J11").Value = 0.2726755566939474
$ws.Range("M11").Value = 0.5599614852851289
$ws.Range("O11").Value = 1.602672744419692
# Row 12
$ws.Range("B12").Value = 1.364502632900098
$ws.Range("C12").Value = 0.3328733956855103
$ws.Range("D12").Value = 0.3711488320999337
$ws.Range("F12").Value = 0.8767278887466858
$ws.Range("G12").Value = 0.3368760154617689
$ws.Range("H12").Value = 0.496200794027402
$ws.Range("I12").Value = 0.3268203491684094
$ws.Range("J12").Value = 0.2724916798820232
$ws.Range("M12").Value = 0.5700587640511259
$ws.Range("O12").Value = 1.600004521615489
# Row 13
$ws.Range("B13").Value = 1.357555465211817
$ws.Range("C13").Value = 0.331191647576901
$ws.Range("D13").Value = 0.3708589359371928
$ws.Range("F13").Value = 0.8767954090823267
$ws.Range("G13").Value = 0.3369165350530423
$ws.Range("H13").Value = 0.4964363947483292
$ws.Range("I13").Value = 0.3272141093171452
$ws.Range("J13").Value = 0.2725300874872048
$ws.Range("M13").Value = 0.5678840351947656
$ws.Range("O13").Value = 1.60056982712463
# Row 14
$ws.Range("B14").Value = 1.334891942424633
$ws.Range("C14").Value = 0.3257049736998852
$ws.Range("D14").Value = 0.3699169764813064
$ws.Range("F14").Value = 0.8770320776903375
$ws.Range("G14").Value = 0.3370579740324331
$ws.Range("H14").Value = 0.4972124431933622
$ws.Range("I14").Value = 0.3285063165772328
$ws.Range("J14").Value = 0.2726598912589253
$ws.Range("M14").Value = 0.5607921489413172
$ws.Range("O14").Value = 1.602449018380455
# Row 15
$ws.Range("B15").Value = 1.321004307831345
$ws.Range("C15").Value = 0.3223426015259747
$ws.Range("D15").Value = 0.3693426559622566
$ws.Range("F15").Value = 0.8771896973479016
$ws.Range("G15").Value = 0.3371517481840058
$ws.Range("H15").Value = 0.4976937139322288
$ws.Range("I15").Value = 0.3293040287907498
$ws.Range("J15").Value = 0.2727428939809187
$ws.Range("M15").Value = 0.5564484641592884
$ws.Range("O15").Value = 1.603627427672279
# Row 16
$ws.Range("B16").Value = 1.241335975414358
$ws.Range("C16").Value = 0.3030496569774925
$ws.Range("D16").Value = 0.3660922622639333
$ws.Range("F16").Value = 0.8782867831701324
$ws.Range("G16").Value = 0.3377983663270356
$ws.Range("H16").Value = 0.5005423703919902
$ws.Range("I16").Value = 0.3339702666277873
$ws.Range("J16").Value = 0.2732720730168623
$ws.Range("M16").Value = 0.5315615984290076
$ws.Range("O16").Value = 1.610799175925891
# Row 17
$ws.Range("B17").Value = 1.192388789855897
$ws.Range("C17").Value = 0.2911925617043494
$ws.Range("D17").Value = 0.3641347520928946
$ws.Range("F17").Value = 0.8791326493958564
$ws.Range("G17").Value = 0.3382923225594965
$ws.Range("H17").Value = 0.5023708206940256
$ws.Range("I17").Value = 0.3369173707628867
$ws.Range("J17").Value = 0.2736444327272665
$ws.Range("M17").Value = 0.5162992769621866
$ws.Range("O17").Value = 1.61557199282295
# Row 18
$ws.Range("B18").Value = 1.164209395604757
$ws.Range("C18").Value = 0.284364935141042
$ws.Range("D18").Value = 0.3630222937393484
$ws.Range("F18").Value = 0.8796826395586024
$ws.Range("G18").Value = 0.338612103836951
$ws.Range("H18").Value = 0.5034522120354339
$ws.Range("I18").Value = 0.3386434727988803
$ws.Range("J18").Value = 0.2738761292627245
$ws.Range("M18").Value = 0.5075228049144869
$ws.Range("O18").Value = 1.618454140335402
# Row 19
$ws.Range("B19").Value = 1.154663887145261
$ws.Range("C19").Value = 0.2820519043078207
$ws.Range("D19").Value = 0.3626479478444082
$ws.Range("F19").Value = 0.8798797533131619
$ws.Range("G19").Value = 0.3387264953202305
$ws.Range("H19").Value = 0.5038234549801146
$ws.Range("I19").Value = 0.3392332213683442
$ws.Range("J19").Value = 0.2739575863569357
$ws.Range("M19").Value = 0.5045516055529333
$ws.Range("O19").Value = 1.619453488754303
# Row 20
$ws.Range("B20").Value = 1.197602034597082
$ws.Range("C20").Value = 0.2924555743468602
$ws.Range("D20").Value = 0.3643417412019403
$ws.Range("F20").Value = 0.8790360354759343
$ws.Range("G20").Value = 0.3382360465217076
$ws.Range("H20").Value = 0.5021731033569665
$ws.Range("I20").Value = 0.3366004362934909
$ws.Range("J20").Value = 0.2736029803068547
$ws.Range("M20").Value = 0.5179237731708639
$ws.Range("O20").Value = 1.615049740184404
# Row 21
$ws.Range("B21").Value = 1.341549933321687
$ws.Range("C21").Value = 0.327316882565384
$ws.Range("D21").Value = 0.3701930993214546
$ws.Range("F21").Value = 0.8769599251852327
$ws.Range("G21").Value = 0.3370149418472579
$ws.Range("H21").Value = 0.4969832654421538
$ws.Range("I21").Value = 0.3281254720311644
$ws.Range("J21").Value = 0.2726210364875641
$ws.Range("M21").Value = 0.5628751462914039
$ws.Range("O21").Value = 1.601891353493642
# Row 22
$ws.Range("B22").Value = 1.435380538648985
$ws.Range("C22").Value = 0.3500284212441898
$ws.Range("D22").Value = 0.3741363102423065
$ws.Range("F22").Value = 0.8761692710014373
$ws.Range("G22").Value = 0.3365361802626552
$ws.Range("H22").Value = 0.4938562782103162
$ws.Range("I22").Value = 0.3228638521634082
$ws.Range("J22").Value = 0.2721356326322919
$ws.Range("M22").Value = 0.5922674117157385
$ws.Range("O22").Value = 1.594515206035481
# Row 23
$ws.Range("B23").Value = 1.385324703346271
$ws.Range("C23").Value = 0.337913630035132
$ws.Range("D23").Value = 0.3720208916006271
$ws.Range("F23").Value = 0.8765394018926145
$ws.Range("G23").Value = 0.3367624084293865
$ws.Range("H23").Value = 0.4955009617031365
$ws.Range("I23").Value = 0.3256466525569746
$ws.Range("J23").Value = 0.2723803812727965
$ws.Range("M23").Value = 0.5765791221098198
$ws.Range("O23").Value = 1.598339836008478
# Row 24
$ws.Range("B24").Value = 1.195245247994649
$ws.Range("C24").Value = 0.2918846000110875
$ws.Range("D24").Value = 0.3642481211107338
$ws.Range("F24").Value = 0.8790795162495115
$ws.Range("G24").Value = 0.3382613774251055
$ws.Range("H24").Value = 0.5022623973249338
$ws.Range("I24").Value = 0.3367436234102961
$ws.Range("J24").Value = 0.2736216660585953
$ws.Range("M24").Value = 0.5171893444974387
$ws.Range("O24").Value = 1.615285420190517
# Row 25
$ws.Range("B25").Value = 0.9892170810735479
$ws.Range("C25").Value = 0.2419408537232357
$ws.Range("D25").Value = 0.3563785493917919
$ws.Range("F25").Value = 0.8842453822170597
$ws.Range("G25").Value = 0.3412409977255635
$ws.Range("H25").Value = 0.510691614558624
$ws.Range("I25").Value = 0.3498982424600108
$ws.Range("J25").Value = 0.2756304262296894
$ws.Range("M25").Value = 0.4532069044750671
$ws.Range("O25").Value = 1.638797315333562
